$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.132.52'
$ws.Range('E2').Value = '  +0.34%  '

$ws.Range('D3').Value = '3.122.90'
$ws.Range('E3').Value = '  +0.59%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.98'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  -0.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.23'
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  +0.43%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = $ws.Range('B7').Style
$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('E8').Value = '  -0.29%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.156'
$ws.Range('D9').Style = $ws.Range('B9').Style
$ws.Range('E9').Value = '  -0.14%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.40'
$ws.Range('D10').Style = $ws.Range('B10').Style
$ws.Range('E10').Value = '  -0.79%  '

$ws.Range('E11').Value = '  -0.80%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000251'
$ws.Range('D12').Style = $ws.Range('B12').Style
$ws.Range('E12').Value = '  +0.60%  '

$ws.Range('E13').Value = '  -0.52%  '

$ws.Range('E14').Value = '  -1.53%  '

$ws.Range('D15').Value = '3.641.02'
$ws.Range('E15').Value = '  +0.66%  '

$ws.Range('D16').Value = '67.103.62'
$ws.Range('E16').Value = '  +0.31%  '

$ws.Range('E17').Value = '  -0.69%  '

$ws.Range('D18').Value = '3.125.79'
$ws.Range('E18').Value = '  +0.65%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.46'
$ws.Range('D19').Style = $ws.Range('B19').Style
$ws.Range('E19').Value = '  +1.82%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '491.44'
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  +1.80%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.94'
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  +5.59%  '

$ws.Range('E22').Value = '  -0.94%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '84.19'
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  +0.24%  '

$ws.Range('E24').Value = '  +0.61%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.30'
$ws.Range('D25').Style = $ws.Range('B25').Style
$ws.Range('E25').Value = '  -3.38%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.36'
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  +2.97%  '

$ws.Range('E27').Value = '  -0.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.90'
$ws.Range('D28').Style = $ws.Range('B28').Style
$ws.Range('E28').Value = '  -0.98%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.35'
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  -1.97%  '

$ws.Range('E30').Value = '  -0.44%  '

$ws.Range('E31').Value = '  -0.12%  '

$ws.Range('E32').Value = '  -0.52%  '

$ws.Range('D33').Value = '0.0₃0953'
$ws.Range('E33').Value = '  -5.61%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('E36').Value = '  -1.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '47.26'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  -1.80%  '

$ws.Range('E38').Value = '  -3.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.310'
$ws.Range('D39').Style = $ws.Range('B39').Style
$ws.Range('E39').Value = '  -2.40%  '

$ws.Range('E40').Value = '  +1.68%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.53'
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  -1.59%  '

$ws.Range('D42').Value = '2.823.88'
$ws.Range('E42').Value = '  -0.28%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.62'
$ws.Range('D43').Style = $ws.Range('B43').Style
$ws.Range('E43').Value = '  -7.46%  '

$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '384.30'
$ws.Range('D44').Style = $ws.Range('B44').Style
$ws.Range('E44').Value = '  -0.23%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0354'
$ws.Range('D45').Style = $ws.Range('B45').Style
$ws.Range('E45').Value = '  -2.43%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.57'
$ws.Range('D46').Style = $ws.Range('B46').Style
$ws.Range('E46').Value = '  +0.72%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.99'
$ws.Range('D48').Style = $ws.Range('B48').Style
$ws.Range('E48').Value = '  +0.38%  '

$ws.Range('E49').Value = '  -1.12%  '

$ws.Range('E50').Value = '  -0.70%  '

$ws.Range('E51').Value = '  -1.12%  '
